$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.046.86"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "1.828.01"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4335"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3680"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07293"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8478"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.74"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("D12").Value = "1.828.10"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.675"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.307"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07068"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.50"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008782"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.94"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").Value = "27.136.77"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.150"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.90"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "2.052.12"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.994"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.48"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.214"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.32"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.248"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.08"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08723"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7434"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.44%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.449"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.910"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.098"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01951"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05250"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.223"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.871"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5132"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.596"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4774"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.17"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.940"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.664"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("E51").Value = "  -1.44%  "
